{"js": "// Insert a new \"DataScience\" paragraph right after the\n// \"Programming Languages: ...\" paragraph (and before \"Microsoft\n// Technologies: ...\").\nconst body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\nlet anchor = null;\nfor (let i = 0; i < paras.items.length; i++) {\n  if (paras.items[i].text.indexOf(\"Programming Languages:\") !== -1) {\n    anchor = paras.items[i];\n    break;\n  }\n}\nif (!anchor) {\n  throw new Error('Anchor paragraph \"Programming Languages:\" not found');\n}\n\nanchor.insertParagraph(\n  \"DataScience :  Machine Learning, Deep Learning, CV,Time Analysis .\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the \"Programming Languages:\" paragraph so the new bullet is\n# inserted right after it (and before \"Microsoft Technologies:\").\n$rng = $d.Content\n$found = $rng.Find.Execute(\"Programming Languages:\")\n$anchorPara = $rng.Paragraphs(1)\n$anchorIndex = $anchorPara.Index\n\n# Insert a new, empty paragraph immediately after the anchor paragraph,\n# then fill it with the new \"DataScience\" line.\n$anchorPara.Range.InsertParagraphAfter()\n$newPara = $d.Paragraphs($anchorIndex + 1)\n$newPara.Range.Text = \"DataScience :  Machine Learning, Deep Learning, CV,Time Analysis .\"\n"}
